# Auto-generated edit script: applies numeric corrections to the
# "currentAveragePrice*" / "LevePrice*" / "LeveProfit*" columns (H:N)
# for specific leve rows across all 8 sheets, per the scheduled-runner diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 383
$ws.Range("I2").Value = 237
$ws.Range("J2").Value = 675
$ws.Range("K2").Value = 237
$ws.Range("L2").Value = 675
$ws.Range("M2").Value = -124
$ws.Range("N2").Value = -901
# Row 6
$ws.Range("H6").Value = 625645.1
$ws.Range("I6").Value = 1428844.6
$ws.Range("J6").Value = 934.44446
$ws.Range("K6").Value = 4286533.800000001
$ws.Range("L6").Value = 2803.33338
$ws.Range("M6").Value = -4286421.800000001
$ws.Range("N6").Value = -3027.33338
# Row 9
$ws.Range("H9").Value = 513.61536
$ws.Range("I9").Value = 591.8889
$ws.Range("J9").Value = 337.5
$ws.Range("K9").Value = 591.8889
$ws.Range("L9").Value = 337.5
$ws.Range("M9").Value = -422.8889
$ws.Range("N9").Value = -675.5
# Row 12
$ws.Range("H12").Value = 24098.285
$ws.Range("I12").Value = 297.6875
$ws.Range("J12").Value = 100260.2
$ws.Range("K12").Value = 297.6875
$ws.Range("L12").Value = 100260.2
$ws.Range("M12").Value = -127.6875
$ws.Range("N12").Value = -100600.2
# Row 15
$ws.Range("H15").Value = 278.44
$ws.Range("I15").Value = 278.44
$ws.Range("K15").Value = 835.3199999999999
$ws.Range("M15").Value = -666.3199999999999
# Row 21
$ws.Range("H21").Value = 12000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 12000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 12000
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -12936
# Row 23
$ws.Range("H23").Value = 12000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 12000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 12000
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -12468
# Row 41
$ws.Range("H41").Value = 1028.7916
$ws.Range("I41").Value = 1447
$ws.Range("K41").Value = 1447
$ws.Range("M41").Value = -1007
# Row 116
$ws.Range("H116").Value = 10033.25
$ws.Range("I116").Value = 10033.25
$ws.Range("K116").Value = 10033.25
$ws.Range("M116").Value = -6591.25
# Row 125
$ws.Range("H125").Value = 6262
$ws.Range("I125").Value = 5516
$ws.Range("J125").Value = 7008
$ws.Range("K125").Value = 49644
$ws.Range("L125").Value = 63072
$ws.Range("M125").Value = -47184
$ws.Range("N125").Value = -67992
# Row 132
$ws.Range("H132").Value = 5562426
$ws.Range("I132").Value = 6586255
$ws.Range("J132").Value = 4497
$ws.Range("K132").Value = 19758765
$ws.Range("L132").Value = 13491
$ws.Range("M132").Value = -19756235
$ws.Range("N132").Value = -18551
# Row 135
$ws.Range("H135").Value = 1536.1708
$ws.Range("I135").Value = 656.37036
$ws.Range("J135").Value = 3232.9285
$ws.Range("K135").Value = 5907.33324
$ws.Range("L135").Value = 29096.3565
$ws.Range("M135").Value = -3372.33324
$ws.Range("N135").Value = -34166.3565
# Row 138
$ws.Range("H138").Value = 3609.2559
$ws.Range("I138").Value = 2641.9285
$ws.Range("J138").Value = 4076.2415
$ws.Range("K138").Value = 7925.7855
$ws.Range("L138").Value = 12228.7245
$ws.Range("M138").Value = -2785.7855
$ws.Range("N138").Value = -22508.7245

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 132
$ws.Range("H132").Value = 16883.205
$ws.Range("I132").Value = 20177.709
$ws.Range("J132").Value = 4117
$ws.Range("K132").Value = 60533.12699999999
$ws.Range("L132").Value = 12351
$ws.Range("M132").Value = -58003.12699999999
$ws.Range("N132").Value = -17411

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 103
$ws.Range("H103").Value = 64500
$ws.Range("J103").Value = 64500
$ws.Range("L103").Value = 64500
$ws.Range("N103").Value = -66844
# Row 134
$ws.Range("H134").Value = 18454.879
$ws.Range("I134").Value = 20822.072
$ws.Range("K134").Value = 62466.216
$ws.Range("M134").Value = -59931.216

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 25
$ws.Range("H25").Value = 15933.333
$ws.Range("I25").Value = 8000
$ws.Range("K25").Value = 8000
$ws.Range("M25").Value = -7826
# Row 58
$ws.Range("H58").Value = 1778.7693
$ws.Range("I58").Value = 1669.6428
$ws.Range("J58").Value = 1906.0834
$ws.Range("K58").Value = 1669.6428
$ws.Range("L58").Value = 1906.0834
$ws.Range("M58").Value = -1466.6428
$ws.Range("N58").Value = -2312.0834
# Row 62
$ws.Range("H62").Value = 2428.6428
$ws.Range("I62").Value = 2000.3334
$ws.Range("K62").Value = 2000.3334
$ws.Range("M62").Value = -1376.3334
# Row 65
$ws.Range("H65").Value = 2428.6428
$ws.Range("I65").Value = 2000.3334
$ws.Range("K65").Value = 10001.667
$ws.Range("M65").Value = -6881.666999999999
# Row 134
$ws.Range("H134").Value = 2093.077
$ws.Range("I134").Value = 928
$ws.Range("J134").Value = 2610.889
$ws.Range("K134").Value = 2784
$ws.Range("L134").Value = 7832.667
$ws.Range("M134").Value = -249
$ws.Range("N134").Value = -12902.667
# Row 136
$ws.Range("H136").Value = 1778.7693
$ws.Range("I136").Value = 1669.6428
$ws.Range("J136").Value = 1906.0834
$ws.Range("K136").Value = 5008.928400000001
$ws.Range("L136").Value = 5718.2502
$ws.Range("M136").Value = -2458.928400000001
$ws.Range("N136").Value = -10818.2502

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 1243.762
$ws.Range("J34").Value = 1304.95
$ws.Range("L34").Value = 3914.85
$ws.Range("N34").Value = -4082.85
# Row 58
$ws.Range("H58").Value = 1800
$ws.Range("J58").Value = 1800
$ws.Range("L58").Value = 5400
$ws.Range("N58").Value = -5656
# Row 122
$ws.Range("H122").Value = 849.6667
$ws.Range("I122").Value = 899
$ws.Range("J122").Value = 825
$ws.Range("K122").Value = 8091
$ws.Range("L122").Value = 7425
$ws.Range("M122").Value = -5641
$ws.Range("N122").Value = -12325
# Row 131
$ws.Range("H131").Value = 852503.9
$ws.Range("I131").Value = 555.9
$ws.Range("J131").Value = 1022893.44
$ws.Range("K131").Value = 1667.7
$ws.Range("L131").Value = 3068680.32
$ws.Range("M131").Value = 3372.3
$ws.Range("N131").Value = -3078760.32

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 52
$ws.Range("H52").Value = 8958.888999999999
$ws.Range("J52").Value = 9575
$ws.Range("L52").Value = 9575
$ws.Range("N52").Value = -10093
# Row 80
$ws.Range("H80").Value = 58885950
$ws.Range("I80").Value = 125130120
$ws.Range("J80").Value = 2246.2222
$ws.Range("K80").Value = 125130120
$ws.Range("L80").Value = 2246.2222
$ws.Range("M80").Value = -125129122
$ws.Range("N80").Value = -4242.2222
# Row 83
$ws.Range("H83").Value = 58885950
$ws.Range("I83").Value = 125130120
$ws.Range("J83").Value = 2246.2222
$ws.Range("K83").Value = 625650600
$ws.Range("L83").Value = 11231.111
$ws.Range("M83").Value = -625645608
$ws.Range("N83").Value = -21215.111
# Row 132
$ws.Range("H132").Value = 3298.7036
$ws.Range("I132").Value = 3033.6667
$ws.Range("J132").Value = 3510.7334
$ws.Range("K132").Value = 9101.000100000001
$ws.Range("L132").Value = 10532.2002
$ws.Range("M132").Value = -6571.000100000001
$ws.Range("N132").Value = -15592.2002

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2925
$ws.Range("I22").Value = 5100
$ws.Range("J22").Value = 750
$ws.Range("K22").Value = 5100
$ws.Range("L22").Value = 750
$ws.Range("M22").Value = -4805
$ws.Range("N22").Value = -1340
# Row 27
$ws.Range("H27").Value = 2925
$ws.Range("I27").Value = 5100
$ws.Range("J27").Value = 750
$ws.Range("K27").Value = 5100
$ws.Range("L27").Value = 750
$ws.Range("M27").Value = -4993
$ws.Range("N27").Value = -964
# Row 46
$ws.Range("H46").Value = 2357.2727
$ws.Range("I46").Value = 628
$ws.Range("J46").Value = 3798.3333
$ws.Range("K46").Value = 628
$ws.Range("L46").Value = 3798.3333
$ws.Range("M46").Value = -440
$ws.Range("N46").Value = -4174.3333
# Row 93
$ws.Range("H93").Value = 1787.4517
$ws.Range("I93").Value = 1688.5
$ws.Range("J93").Value = 2029.3334
$ws.Range("K93").Value = 1688.5
$ws.Range("L93").Value = 2029.3334
$ws.Range("M93").Value = -440.5
$ws.Range("N93").Value = -4525.3334
# Row 106
$ws.Range("H106").Value = 37116
$ws.Range("J106").Value = 37116
$ws.Range("L106").Value = 37116
$ws.Range("M106").Value = -39640

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2835.7646
$ws.Range("I122").Value = 2015.6154
$ws.Range("K122").Value = 6046.8462
$ws.Range("M122").Value = -3596.8462
